# The "Credit" column (C) text for the three hymn/song rows is removed,
# leaving those cells blank (their style is kept). This also makes the
# three long "credit" strings in the shared-strings table unused, and the
# narrow "Credit" column no longer needs to be as wide, so the column is
# shrunk back down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the credit text out of C3 (Be Thou My Vision), C5 (You Raise Me Up)
# and C6 (10,000 Reasons) - formatting/style stays, only the value goes away.
$ws.Range("C3").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("C6").ClearContents()

# Column C no longer needs to hold the long credit paragraphs, shrink it
# back down to a normal width.
$ws.Columns.Item(3).ColumnWidth = 17

# Move the active selection to C3 (where the author had been working).
$ws.Range("C3").Select()
